# Script 1 - atualização automática de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: "Variação (%) 2022" -> "Variação (%) 2023" (rows 2-10)
$ws.Range("A2").Value = "AC"
$ws.Range("B2").Value = 14.73257689442189
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Variação (%) 2023"

$ws.Range("A3").Value = "MS"
$ws.Range("B3").Value = 13.44269577606423
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "Variação (%) 2023"

$ws.Range("A4").Value = "MT"
$ws.Range("B4").Value = 12.88001598426398
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "Variação (%) 2023"

$ws.Range("A5").Value = "TO"
$ws.Range("B5").Value = 7.890383025089162
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "Variação (%) 2023"

$ws.Range("A6").Value = "RJ"
$ws.Range("B6").Value = 5.652659822157795
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Variação (%) 2023"

$ws.Range("A7").Value = "GO"
$ws.Range("B7").Value = 4.816953216278661
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = "Variação (%) 2023"

$ws.Range("A8").Value = "SE"
$ws.Range("B8").Value = 3.118144130554446
$ws.Range("C8").Value = 15
$ws.Range("D8").Value = "Variação (%) 2023"

$ws.Range("A9").Value = "BR"
$ws.Range("B9").Value = 3.241657824791806
$ws.Range("D9").Value = "Variação (%) 2023"

$ws.Range("A10").Value = "NE"
$ws.Range("B10").Value = 2.867008788862638
$ws.Range("D10").Value = "Variação (%) 2023"

# Block 2: "Variação (%) 2022/2010" -> "Variação (%) 2023/2010" (rows 11-19)
$ws.Range("A11").Value = "MT"
$ws.Range("B11").Value = 77.14346626765018
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Variação (%) 2023/2010"

$ws.Range("A12").Value = "RR"
$ws.Range("B12").Value = 63.84407511155798
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = "Variação (%) 2023/2010"

$ws.Range("A13").Value = "TO"
$ws.Range("B13").Value = 59.35173933449352
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = "Variação (%) 2023/2010"

$ws.Range("A14").Value = "MS"
$ws.Range("B14").Value = 49.57179111911111
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = "Variação (%) 2023/2010"

$ws.Range("A15").Value = "AC"
$ws.Range("B15").Value = 42.69524774665621
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = "Variação (%) 2023/2010"

$ws.Range("A16").Value = "PI"
$ws.Range("B16").Value = 36.70681689547283
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = "Variação (%) 2023/2010"

$ws.Range("A17").Value = "SE"
$ws.Range("B17").Value = 7.324239245718005
$ws.Range("C17").Value = 25
$ws.Range("D17").Value = "Variação (%) 2023/2010"

$ws.Range("A18").Value = "BR"
$ws.Range("B18").Value = 15.14430956101356
$ws.Range("D18").Value = "Variação (%) 2023/2010"

$ws.Range("A19").Value = "NE"
$ws.Range("B19").Value = 16.86384673819174
$ws.Range("D19").Value = "Variação (%) 2023/2010"
